# Update the date embedded in the worksheet's status text.
# (source: "Date: 05-10-2018 - Department: Sales department"
#  target: "Date: 06-01-2026 - Department: Sales department")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Date: 06-01-2026 - Department: Sales department"
